$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"): copy formatting from the
# existing header cell H1 so they share the same bold/bordered/centered style.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for new columns I and J
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 5

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 4

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 6

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2
